# Spring 2025 Presentation - "Added App dev 2 part"
#
# Inserts a new slide ("App Dev 2") as the 2nd slide in the deck (pushing
# the existing "Design 4 All" and "Core- Communications" slides down one
# position), using the same "Title and Content" layout as those slides.

$p = $ppt.ActivePresentation

# The master's 2nd custom layout is "Title and Content" - the same layout
# used by the existing "Design 4 All" / "Core- Communications" slides.
$layout = $p.SlideMaster.CustomLayouts.Item(2)

# Insert the new slide at position 2 (after the title slide).
$slide = $p.Slides.AddSlide(2, $layout)

# --- Title -----------------------------------------------------------
$slide.Shapes.Item(1).TextFrame.TextRange.Text = "App Dev 2  "

# --- Body content ------------------------------------------------------
$body = $slide.Shapes.Item(2)
$tr = $body.TextFrame.TextRange

$tr.Text = "Worked on an address book"
[void]$tr.InsertAfter("`rLearned more about classes")
[void]$tr.InsertAfter("`rFigured out how to put code into classes, thereby allowing the organization of code within multiple files")
[void]$tr.InsertAfter("`rRealized even moreso how bad I am at UI design...")

# Demote the two "classes" bullets to a sub-level with a Courier New "o" bullet.
$p2 = $tr.Paragraphs(2)
$p2.IndentLevel = 2
$p2.ParagraphFormat.Bullet.Font.Name = "Courier New"
$p2.ParagraphFormat.Bullet.Character = 111

$p3 = $tr.Paragraphs(3)
$p3.IndentLevel = 2
$p3.ParagraphFormat.Bullet.Font.Name = "Courier New"
$p3.ParagraphFormat.Bullet.Character = 111

# Match the body placeholder formatting used on the other content slides.
$body.TextFrame.MarginLeft = 7.2
$body.TextFrame.MarginRight = 7.2
$body.TextFrame.MarginTop = 3.6
$body.TextFrame.MarginBottom = 3.6
$body.TextFrame.Orientation = 1
$body.TextFrame2.VerticalAnchor = 1
$body.TextFrame.AutoSize = 2
